$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -58
$ws.Range("B3").Value = -89.40000000000001
$ws.Range("C3").Value = 58
$ws.Range("C4").Value = -35
$ws.Range("C5").Value = 9.9
$ws.Range("C6").Value = 25
$ws.Range("C7").Value = 48
$ws.Range("C8").Value = 71
$ws.Range("C9").Value = -20.2
$ws.Range("C10").Value = -33.8
$ws.Range("C11").Value = -45.8
$ws.Range("C12").Value = -33.9
$ws.Range("C13").Value = -5.9
$ws.Range("C14").Value = 58.6
$ws.Range("C15").Value = 124.6
$ws.Range("C16").Value = 100.4
$ws.Range("C17").Value = 109.3
$ws.Range("C18").Value = 84.8
$ws.Range("C19").Value = 40.7
$ws.Range("C20").Value = 85.09999999999999
$ws.Range("C22").Value = 108.3
$ws.Range("C23").Value = 92.90000000000001
$ws.Range("C24").Value = 113.3
